# The sentence "...هیچ دستور پخت غذایی در فهرست دستور پخت ها وجود داشته باشد."
# needs to become "...هیچ دستور پخت غذایی در فهرست دستور پخت ها وجود نداشته باشد."
# (the recipe-exists assertion is negated -> "recipe removed" check, issue #26).
#
# Only that single paragraph should change, so we locate it first by its
# unique wording and then run Find/Replace scoped to just that paragraph's
# Range (never against the whole document) to avoid touching the many other
# "... وجود داشته باشد" occurrences elsewhere in the doc.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*هیچ دستور پخت غذایی در فهرست دستور پخت ها وجود داشته باشد*") {
        $p.Range.Find.Execute(
            " وجود داشته باشد",  # Find what
            $true,                # MatchCase
            $false,               # MatchWholeWord
            $false,               # MatchWildcards
            $false,               # MatchSoundsLike
            $false,               # MatchAllWordForms
            $true,                # Forward
            1,                    # Wrap (wdFindContinue)
            $false,               # Format
            " وجود نداشته باشد", # Replace with
            2                     # Replace (wdReplaceOne)
        )
        break
    }
}
